{"js": "// 1) Refresh the \"retrieved\" timestamp in the footer.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst primarySection = sections.items[0];\nconst footer = primarySection.getFooter(Word.HeaderFooterType.primary);\nconst searchResults = footer.search(\"2025-06-30 12:12Z\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  searchResults.items[0].insertText(\"2025-07-02 02:48Z\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Add the standard set of inline-markup character styles used by the\n//    PubMed -> Word converter (b/i/sub/sup/u), each based on the built-in\n//    \"Default Paragraph Font\" style.\nconst styleDefs = [\n  { id: \"b\", apply: (f) => { f.bold = true; } },\n  { id: \"i\", apply: (f) => { f.italic = true; } },\n  { id: \"sub\", apply: (f) => { f.subscript = true; } },\n  { id: \"sup\", apply: (f) => { f.superscript = true; } },\n  { id: \"u\", apply: (f) => { f.underline = Word.UnderlineType.single; } }\n];\n\nfor (const def of styleDefs) {\n  context.document.addStyle(def.id, Word.StyleType.character);\n  await context.sync();\n\n  const style = context.document.styles.getByNameOrNullObject(def.id);\n  style.baseStyle = \"DefaultParagraphFont\";\n  style.priority = 1;\n  style.quickStyle = true;\n  def.apply(style.font);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Refresh the \"retrieved\" timestamp in the footer.\n$footer = $d.Sections(1).Footers(1)\n$find = $footer.Range.Find\n$find.Text = \"2025-06-30 12:12Z\"\n$find.Replacement.Text = \"2025-07-02 02:48Z\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# 2) Add the standard set of inline-markup character styles used by the\n#    PubMed -> Word converter (b/i/sub/sup/u), each based on the built-in\n#    \"Default Paragraph Font\" style.\n$sb = $d.Styles.Add(\"b\", 2)\n$sb.BaseStyle = \"DefaultParagraphFont\"\n$sb.Priority = 1\n$sb.QuickStyle = 1\n$sb.Font.Bold = 1\n\n$si = $d.Styles.Add(\"i\", 2)\n$si.BaseStyle = \"DefaultParagraphFont\"\n$si.Priority = 1\n$si.QuickStyle = 1\n$si.Font.Italic = 1\n\n$ssub = $d.Styles.Add(\"sub\", 2)\n$ssub.BaseStyle = \"DefaultParagraphFont\"\n$ssub.Priority = 1\n$ssub.QuickStyle = 1\n$ssub.Font.Subscript = $true\n\n$ssup = $d.Styles.Add(\"sup\", 2)\n$ssup.BaseStyle = \"DefaultParagraphFont\"\n$ssup.Priority = 1\n$ssup.QuickStyle = 1\n$ssup.Font.Superscript = $true\n\n$su = $d.Styles.Add(\"u\", 2)\n$su.BaseStyle = \"DefaultParagraphFont\"\n$su.Priority = 1\n$su.QuickStyle = 1\n$su.Font.Underline = 1\n"}
